$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 819
$ws.Range("F3").Value = 14741
$ws.Range("F4").Value = 16
$ws.Range("F5").Value = 1658
$ws.Range("F6").Value = 504
$ws.Range("F8").Value = 1305
$ws.Range("F9").Value = 1988
$ws.Range("F10").Value = 947
$ws.Range("F11").Value = 44
$ws.Range("F12").Value = 2350
$ws.Range("F13").Value = 617
$ws.Range("F14").Value = 840
$ws.Range("F15").Value = 3660
$ws.Range("F18").Value = 2709
$ws.Range("F19").Value = 695
$ws.Range("F20").Value = 131
$ws.Range("F22").Value = 1927
$ws.Range("F23").Value = 1139
$ws.Range("F24").Value = 1641
$ws.Range("F25").Value = 346
$ws.Range("F26").Value = 178
$ws.Range("F27").Value = 7597
$ws.Range("F28").Value = 5258
$ws.Range("F29").Value = 333
$ws.Range("F31").Value = 729
$ws.Range("F32").Value = 733
$ws.Range("F33").Value = 3409
$ws.Range("F36").Value = 362
$ws.Range("F38").Value = 124
$ws.Range("F39").Value = 4517
$ws.Range("F40").Value = 742
$ws.Range("F41").Value = 33
$ws.Range("F42").Value = 350

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 11
$ws.Range("F15").Value = 100
$ws.Range("F18").Value = 121
$ws.Range("F19").Value = 60

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8051
$ws.Range("F3").Value = 319
$ws.Range("F4").Value = 1143

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 8051
$ws.Range("F3").Value = 819
$ws.Range("F4").Value = 319
$ws.Range("F5").Value = 1143
$ws.Range("F6").Value = 14741
$ws.Range("F9").Value = 1658
$ws.Range("F10").Value = 504
$ws.Range("F11").Value = 1305
$ws.Range("F12").Value = 1988
$ws.Range("F14").Value = 44
$ws.Range("F15").Value = 617
$ws.Range("F17").Value = 3660
$ws.Range("F19").Value = 2709
$ws.Range("F20").Value = 696
$ws.Range("F21").Value = 131
$ws.Range("F23").Value = 1927
$ws.Range("F24").Value = 11
$ws.Range("F29").Value = 1641
$ws.Range("F30").Value = 100
$ws.Range("F31").Value = 346
$ws.Range("F32").Value = 178
$ws.Range("F33").Value = 7599
$ws.Range("F34").Value = 5258
$ws.Range("F35").Value = 333
$ws.Range("F36").Value = 729
$ws.Range("F37").Value = 3409
$ws.Range("F40").Value = 362
$ws.Range("F42").Value = 124
$ws.Range("F43").Value = 4517
$ws.Range("F44").Value = 742
$ws.Range("F45").Value = 33
$ws.Range("F46").Value = 350
